# Generate Report for Handoff
# Row 3 (7e5ed80a-e98d-40f5-9011-06681c95db56.md) is now "Ready for handoff",
# and the "Latest Handoff Datetime" stamps for both rows are refreshed to
# reflect the newly generated handoff report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest Handoff Date (col D) refreshed for both rows,
# and row 3's status (zh-cn/de-de columns B & C) becomes "Ready for handoff".
$wsOverview.Range("D2").Value = "2016-03-24 11:00:40"
$wsOverview.Range("D3").Value = "2016-03-24 11:00:40"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: Latest Handoff Datetime (col E) refreshed for both rows,
# and row 3's Status (col C) becomes "Ready for handoff".
$wsZhCn.Range("E2").Value = "2016-03-24 11:00:26"
$wsZhCn.Range("E3").Value = "2016-03-24 11:00:26"
$wsZhCn.Range("C3").Value = "Ready for handoff"

# de-de sheet: Latest Handoff Datetime (col E) refreshed for both rows,
# and row 3's Status (col C) becomes "Ready for handoff".
$wsDeDe.Range("E2").Value = "2016-03-24 11:00:40"
$wsDeDe.Range("E3").Value = "2016-03-24 11:00:40"
$wsDeDe.Range("C3").Value = "Ready for handoff"
